$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entries to insert, in ascending order of their FINAL row position.
# Each row: final row number, Name (column A), hex code (column D) -- the
# PUA character for column B and its decimal code point for column C are
# derived from the hex code.
$entries = @(
  @{ Row = 1;   Name = "Abel";                       Hex = "ea96" },
  @{ Row = 18;  Name = "Barachiah";                   Hex = "ea97" },
  @{ Row = 31;  Name = "Chorazin";                    Hex = "ea8f" },
  @{ Row = 53;  Name = "Herodians";                   Hex = "ea05" },
  @{ Row = 60;  Name = "Israel";                      Hex = "ea04" },
  @{ Row = 69;  Name = "Jericho";                     Hex = "ea99" },
  @{ Row = 76;  Name = "Jonah(Simon" + [char]0x2019 + "s_father)"; Hex = "ea29" },
  @{ Row = 90;  Name = "Magadan";                     Hex = "ea30" },
  @{ Row = 130; Name = "Simon(Jesus" + [char]0x2019 + "_Brother)"; Hex = "ea92" },
  @{ Row = 134; Name = "Sodom";                       Hex = "ea8e" },
  @{ Row = 149; Name = "Zechariah";                   Hex = "ea98" }
)

foreach ($entry in $entries) {
  $r = $entry.Row
  $ws.Rows("$($r):$($r)").Insert()

  $dec = [Convert]::ToInt32($entry.Hex, 16)
  $puaChar = [char]$dec

  $ws.Cells.Item($r, 1).Value = $entry.Name
  $ws.Cells.Item($r, 2).Value = $puaChar
  $ws.Cells.Item($r, 3).Value = $dec
  $ws.Cells.Item($r, 4).Value = $entry.Hex
}

# Restore the view the workbook shipped with after the edit: top-left A1,
# active cell D1 (matches the post-edit sheetView in the target workbook).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("D1").Select()
